$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New handed-back file being reported on: 436f300c-1f5e-491c-80d3-ad8b05589f64.md
# A new row (row 4) is appended to each of the three worksheets.
# ---------------------------------------------------------------------------

$newMd   = "436f300c-1f5e-491c-80d3-ad8b05589f64.md"
$zhXlf   = "436f300c-1f5e-491c-80d3-ad8b05589f64.782b63dff3d1b5d2675e0b3518f26e9027adfde9.zh-cn.xlf"
$deXlf   = "436f300c-1f5e-491c-80d3-ad8b05589f64.782b63dff3d1b5d2675e0b3518f26e9027adfde9.de-de.xlf"
$inSync  = "Handed back: in sync with en-US"

# =====================  Sheet 1: "Overview"  ===============================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $newMd
$wsOverview.Range("B4").Value = $inSync
$wsOverview.Range("C4").Value = $inSync

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/436f300c1f5e491c80d3ad8b05589f64commit/e2e/436f300c-1f5e-491c-80d3-ad8b05589f64.md",
    "",
    "",
    $newMd
)

# =====================  Sheet 2: "zh-cn"  ===================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $newMd
$wsZh.Range("B4").Value = $inSync
$wsZh.Range("C4").Value = $zhXlf
$wsZh.Range("D4").Value = "2016-02-19 07:52:33"
$wsZh.Range("E4").Value = $newMd
$wsZh.Range("F4").Value = $zhXlf
$wsZh.Range("G4").Value = "2016-02-19 07:53:29"
$wsZh.Range("H4").Value = "Include"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/436f300c1f5e491c80d3ad8b05589f64commit/e2e/436f300c-1f5e-491c-80d3-ad8b05589f64.md",
    "",
    "",
    $newMd
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/436f300c1f5e491c80d3ad8b05589f64commit/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/436f300c-1f5e-491c-80d3-ad8b05589f64.782b63dff3d1b5d2675e0b3518f26e9027adfde9.zh-cn.xlf",
    "",
    "",
    $zhXlf
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/436f300c1f5e491c80d3ad8b05589f64commit/e2e/436f300c-1f5e-491c-80d3-ad8b05589f64.md",
    "",
    "",
    $newMd
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/436f300c1f5e491c80d3ad8b05589f64commit/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/436f300c-1f5e-491c-80d3-ad8b05589f64.782b63dff3d1b5d2675e0b3518f26e9027adfde9.zh-cn.xlf",
    "",
    "",
    $zhXlf
)

# =====================  Sheet 3: "de-de"  ===================================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $newMd
$wsDe.Range("B4").Value = $inSync
$wsDe.Range("C4").Value = $deXlf
$wsDe.Range("D4").Value = "2016-02-19 07:52:46"
$wsDe.Range("E4").Value = $newMd
$wsDe.Range("F4").Value = $deXlf
$wsDe.Range("G4").Value = "2016-02-19 07:53:54"
$wsDe.Range("H4").Value = "Include"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/436f300c1f5e491c80d3ad8b05589f64commit/e2e/436f300c-1f5e-491c-80d3-ad8b05589f64.md",
    "",
    "",
    $newMd
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/436f300c1f5e491c80d3ad8b05589f64commit/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/436f300c-1f5e-491c-80d3-ad8b05589f64.782b63dff3d1b5d2675e0b3518f26e9027adfde9.de-de.xlf",
    "",
    "",
    $deXlf
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/436f300c1f5e491c80d3ad8b05589f64commit/e2e/436f300c-1f5e-491c-80d3-ad8b05589f64.md",
    "",
    "",
    $newMd
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/436f300c1f5e491c80d3ad8b05589f64commit/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/436f300c-1f5e-491c-80d3-ad8b05589f64.782b63dff3d1b5d2675e0b3518f26e9027adfde9.de-de.xlf",
    "",
    "",
    $deXlf
)

Write-Host "Report row added for 436f300c-1f5e-491c-80d3-ad8b05589f64.md on all sheets."
